# "Adjust animation and move balance"
#
# frameData.xlsx (sheet "Arkusz1") keeps a mirrored animation-frame table:
# rows 3-23 simply re-read rows 27-47 via formulas, so editing the raw
# numbers in the lower table (row 29 = "Standing Right Kick") ripples back
# up into the matching top-table row (row 5) automatically.
#
# The hit/knockback balance for that frame is adjusted:
#   B29 (hit stun)      20 -> 18
#   D29 (knockback)    -12 -> -14
# F29/G29 (and the mirrored B5/D5/F5/G5) are formulas and recompute on
# their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B29").Value = 18
$ws.Range("D29").Value = -14

# Move the selection/cursor down to the balance row that was just edited.
$ws.Range("B30").Select()
